$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "asad"
$ws.Range("B3").Value = "hghghg"

$ws.Range("B3").Select()
